$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.266.26"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").Value = "2.523.38"
$ws.Range("E3").Value = "  +1.09%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.07"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.06"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.31%  "

# Row 7
$ws.Range("E7").Value = "  +0.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.68%  "

# Row 9
$ws.Range("D9").Value = "2.521.48"
$ws.Range("E9").Value = "  +0.98%  "

# Row 10
$ws.Range("E10").Value = "  -1.14%  "

# Row 11
$ws.Range("E11").Value = "  -1.53%  "

# Row 12
$ws.Range("E12").Value = "  -3.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "

# Row 14
$ws.Range("D14").Value = "2.966.25"
$ws.Range("E14").Value = "  +1.01%  "

# Row 15
$ws.Range("D15").Value = "58.304.34"
$ws.Range("E15").Value = "  +0.14%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.13"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("D18").Value = "2.519.48"
$ws.Range("E18").Value = "  +0.93%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.66"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "321.89"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.80%  "

# Row 22
$ws.Range("E22").Value = "  +7.06%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.48"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

# Row 25
$ws.Range("E25").Value = "  -1.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.07%  "

# Row 28
$ws.Range("E28").Value = "  +0.13%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0751"
$ws.Range("E29").Value = "  +0.21%  "

# Row 30
$ws.Range("E30").Value = "  -0.68%  "

# Row 31
$ws.Range("E31").Value = "  +1.38%  "

# Row 32
$ws.Range("E32").Value = "  -0.22%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.17%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.16"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "

# Row 37
$ws.Range("E37").Value = "  -6.39%  "

# Row 38
$ws.Range("E38").Value = "  -2.83%  "

# Row 39
$ws.Range("E39").Value = "  +0.92%  "

# Row 40
$ws.Range("E40").Value = "  -0.60%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.769"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.61%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "276.41"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.38%  "

# Row 43
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "129.99"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.85%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.98"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.12%  "

# Row 46
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("E47").Value = "  +0.81%  "

# Row 48
$ws.Range("E48").Value = "  +1.84%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.69"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.63%  "

# Row 50
$ws.Range("E50").Value = "  +0.01%  "

# Row 51
$ws.Range("E51").Value = "  -0.83%  "
